$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 7 and 8 with revised values
$ws.Range("A7").Value = 75
$ws.Range("B7").Value = 342
$ws.Range("A8").Value = 80
$ws.Range("B8").Value = 387

# Add new empty, but styled, "marker" cells below the data (rows 11 and 14)
$ws.Range("A11").NumberFormat = "#,##0"
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").Font.Underline = 2

$ws.Range("B14").NumberFormat = "#,##0"
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("B14").Font.Underline = 2

# Move the active selection to B9
$ws.Range("B9").Select()
